# Blackjack strategy workbook: add a new "split" strategy sheet and
# duplicate the text-highlight conditional formatting on "hard" that was
# (re-)applied while the new sheet's formatting was being set up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "split" sheet after the existing "soft" sheet.
# ---------------------------------------------------------------------
$hard = $wb.Worksheets.Item("hard")
$soft = $wb.Worksheets.Item("soft")

$split = $wb.Worksheets.Add([System.Type]::Missing, $soft)
$split.Name = "split"

# ---------------------------------------------------------------------
# 2. Populate the "split" sheet: header row + pair-value rows 2..11,
#    every strategy cell defaulting to "No Split".
# ---------------------------------------------------------------------
$headers = @("Player","Dealer2","Dealer3","Dealer4","Dealer5","Dealer6","Dealer7","Dealer8","Dealer9","Dealer10","Dealer11")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $split.Cells.Item(1, $i + 1).Value = $headers[$i]
}

for ($r = 2; $r -le 11; $r++) {
    $split.Cells.Item($r, 1).Value = $r
    for ($c = 2; $c -le 11; $c++) {
        $split.Cells.Item($r, $c).Value = "No Split"
    }
}

# ---------------------------------------------------------------------
# 3. Conditional formatting on "split": dealer 7-11 (G:K) and dealer 2-6
#    (B:F) are two separate rule sets, each highlighting "No Split" in
#    grey and "Surrender" in white.
# ---------------------------------------------------------------------
$greyColor = 14277081   # D9D9D9 - White, Background 1, Darker 15%
$whiteColor = 16777215  # FFFFFF - White, Background 1

$rngStrong = $split.Range("G2:K11")
$condNoSplitStrong = $rngStrong.FormatConditions.Add(2, 3, '=NOT(ISERROR(SEARCH("No Split",G2)))')
$condNoSplitStrong.Interior.Color = $greyColor
$condSurrenderStrong = $rngStrong.FormatConditions.Add(2, 3, '=NOT(ISERROR(SEARCH("Surrender",G2)))')
$condSurrenderStrong.Interior.Color = $whiteColor

$rngWeak = $split.Range("B2:F11")
$condNoSplitWeak = $rngWeak.FormatConditions.Add(2, 3, '=NOT(ISERROR(SEARCH("No Split",B2)))')
$condNoSplitWeak.Interior.Color = $greyColor
$condSurrenderWeak = $rngWeak.FormatConditions.Add(2, 3, '=NOT(ISERROR(SEARCH("Surrender",B2)))')
$condSurrenderWeak.Interior.Color = $whiteColor

# ---------------------------------------------------------------------
# 4. "hard" picked up a second copy of its text-highlight conditional
#    formatting rules (Split/Double/Stand/Hit) over the same B2:K19
#    range, on top of the pre-existing set.
# ---------------------------------------------------------------------
$splitColor = 1137094    # C65911 - Orange, Accent 2, Darker 25%
$doubleColor = 10086143  # FFE699 - Gold, Accent 4, Lighter 60%
$standColor = 14395790   # 8EA9DB - Blue, Accent 1, Lighter 40%
$hitColor = 11854022     # C6E0B4 - Green, Accent 6, Lighter 60%

$rngHard = $hard.Range("B2:K19")

$condSplit = $rngHard.FormatConditions.Add(2, 3, '=NOT(ISERROR(SEARCH("Split",B2)))')
$condSplit.Interior.Color = $splitColor

$condDouble = $rngHard.FormatConditions.Add(2, 3, '=NOT(ISERROR(SEARCH("Double",B2)))')
$condDouble.Interior.Color = $doubleColor

$condStand = $rngHard.FormatConditions.Add(2, 3, '=NOT(ISERROR(SEARCH("Stand",B2)))')
$condStand.Interior.Color = $standColor

$condHit = $rngHard.FormatConditions.Add(2, 3, '=NOT(ISERROR(SEARCH("Hit",B2)))')
$condHit.Interior.Color = $hitColor

# ---------------------------------------------------------------------
# 5. Make "split" the active sheet / selection, matching where the user
#    left off editing.
# ---------------------------------------------------------------------
$split.Activate()
$split.Range("L9").Select()

Write-Host "Added 'split' sheet and refreshed conditional formatting"
